# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice/LevePrice/LeveProfit figures (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets with
# newly pulled market-board data.

$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 50500
$ws.Range("J17").Value = 50500
$ws.Range("L17").Value = 151500
$ws.Range("N17").Value = -151836
$ws.Range("H33").Value = 236.5
$ws.Range("I33").Value = 236.5
$ws.Range("K33").Value = 236.5
$ws.Range("M33").Value = -7.5
$ws.Range("H43").Value = 5260.2
$ws.Range("J43").Value = 4266.6665
$ws.Range("L43").Value = 4266.6665
$ws.Range("N43").Value = -4404.6665
$ws.Range("H98").Value = 937.7
$ws.Range("I98").Value = 930.8889
$ws.Range("K98").Value = 930.8889
$ws.Range("M98").Value = 567.1111
$ws.Range("H122").Value = 937.7
$ws.Range("I122").Value = 930.8889
$ws.Range("K122").Value = 2792.6667
$ws.Range("M122").Value = -342.6667000000002
$ws.Range("H125").Value = 4500
$ws.Range("I125").Value = 4500
$ws.Range("K125").Value = 40500
$ws.Range("M125").Value = -38040
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1513.4286
$ws.Range("I5").Value = 1879.6
$ws.Range("J5").Value = 598
$ws.Range("K5").Value = 1879.6
$ws.Range("L5").Value = 598
$ws.Range("M5").Value = -1767.6
$ws.Range("N5").Value = -822
$ws.Range("H21").Value = 4001.875
$ws.Range("I21").Value = 419.16666
$ws.Range("J21").Value = 14750
$ws.Range("K21").Value = 419.16666
$ws.Range("L21").Value = 14750
$ws.Range("M21").Value = -45.16665999999998
$ws.Range("N21").Value = -15498
$ws.Range("H32").Value = 8562091
$ws.Range("I32").Value = 8757353
$ws.Range("K32").Value = 8757353
$ws.Range("M32").Value = -8757066
$ws.Range("H45").Value = 4498.909
$ws.Range("I45").Value = 4311
$ws.Range("K45").Value = 4311
$ws.Range("M45").Value = -3934
$ws.Range("H113").Value = 149990.5
$ws.Range("J113").Value = 149990.5
$ws.Range("L113").Value = 149990.5
$ws.Range("N113").Value = -158668.5
$ws.Range("H122").Value = 22110.25
$ws.Range("I122").Value = 24697.428
$ws.Range("K122").Value = 74092.284
$ws.Range("M122").Value = -71642.284

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1513.4286
$ws.Range("I4").Value = 1879.6
$ws.Range("J4").Value = 598
$ws.Range("K4").Value = 1879.6
$ws.Range("L4").Value = 598
$ws.Range("M4").Value = -1764.6
$ws.Range("N4").Value = -828
$ws.Range("H86").Value = 3532
$ws.Range("J86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 3532
$ws.Range("J89").Value = 5000
$ws.Range("L89").Value = 25000
$ws.Range("N89").Value = -36232
$ws.Range("H94").Value = 1254.4828
$ws.Range("I94").Value = 1410.9412
$ws.Range("K94").Value = 1410.9412
$ws.Range("M94").Value = -959.9412
$ws.Range("H105").Value = 2579.25
$ws.Range("I105").Value = 2579.25
$ws.Range("K105").Value = 2579.25
$ws.Range("M105").Value = -832.25

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1700
$ws.Range("H58").Value = 2113.0557
$ws.Range("I58").Value = 1405.125
$ws.Range("J58").Value = 7776.5
$ws.Range("K58").Value = 1405.125
$ws.Range("L58").Value = 7776.5
$ws.Range("M58").Value = -1202.125
$ws.Range("N58").Value = -8182.5
$ws.Range("H136").Value = 2113.0557
$ws.Range("I136").Value = 1405.125
$ws.Range("J136").Value = 7776.5
$ws.Range("K136").Value = 4215.375
$ws.Range("L136").Value = 23329.5
$ws.Range("M136").Value = -1665.375
$ws.Range("N136").Value = -28429.5

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 118.2
$ws.Range("I2").Value = 37.142857
$ws.Range("J2").Value = 189.125
$ws.Range("K2").Value = 222.857142
$ws.Range("L2").Value = 1134.75
$ws.Range("M2").Value = -109.857142
$ws.Range("N2").Value = -1360.75
$ws.Range("H38").Value = 2099
$ws.Range("I38").Value = 4208
$ws.Range("J38").Value = 411.8
$ws.Range("K38").Value = 12624
$ws.Range("L38").Value = 1235.4
$ws.Range("M38").Value = -12277
$ws.Range("N38").Value = -1929.4
$ws.Range("H94").Value = 16200
$ws.Range("J94").Value = 16200
$ws.Range("L94").Value = 48600
$ws.Range("N94").Value = -49952
$ws.Range("H113").Value = 1533
$ws.Range("J113").Value = 1533
$ws.Range("L113").Value = 4599
$ws.Range("N113").Value = -8939
$ws.Range("H125").Value = 6600
$ws.Range("J125").Value = 6600
$ws.Range("L125").Value = 19800
$ws.Range("N125").Value = -29640
$ws.Range("H137").Value = 2985.889
$ws.Range("J137").Value = 3403.1667
$ws.Range("L137").Value = 10209.5001
$ws.Range("N137").Value = -20409.5001

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1500
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -2492
$ws.Range("H102").Value = 2251.0715
$ws.Range("I102").Value = 1990.8889
$ws.Range("K102").Value = 1990.8889
$ws.Range("M102").Value = -368.8888999999999

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 809
$ws.Range("I22").Value = 699.75
$ws.Range("K22").Value = 699.75
$ws.Range("M22").Value = -404.75
$ws.Range("H27").Value = 809
$ws.Range("I27").Value = 699.75
$ws.Range("K27").Value = 699.75
$ws.Range("M27").Value = -592.75
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""
$ws.Range("H46").Value = 1612.16
$ws.Range("I46").Value = 1257.75
$ws.Range("J46").Value = 3029.8
$ws.Range("K46").Value = 1257.75
$ws.Range("L46").Value = 3029.8
$ws.Range("M46").Value = -1069.75
$ws.Range("N46").Value = -3405.8
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H61").Value = 1499
$ws.Range("I61").Value = 1499
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1499
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1297
$ws.Range("N61").Value = ""
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352
$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 671
$ws.Range("N113").Value = ""
$ws.Range("H120").Value = 35560
$ws.Range("J120").Value = 35560
$ws.Range("L120").Value = 35560
$ws.Range("N120").Value = -45236
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").Value = ""

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 18952.4
$ws.Range("I74").Value = 18777
$ws.Range("J74").Value = 18996.25
$ws.Range("K74").Value = 18777
$ws.Range("L74").Value = 18996.25
$ws.Range("M74").Value = -17841
$ws.Range("N74").Value = -20868.25
$ws.Range("H77").Value = 18952.4
$ws.Range("I77").Value = 18777
$ws.Range("J77").Value = 18996.25
$ws.Range("K77").Value = 56331
$ws.Range("L77").Value = 56988.75
$ws.Range("M77").Value = -51651
$ws.Range("N77").Value = -66348.75
$ws.Range("H81").Value = 8489.637000000001
$ws.Range("I81").Value = 8028.5713
$ws.Range("K81").Value = 16057.1426
$ws.Range("M81").Value = -14996.1426
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = ""
$ws.Range("H84").Value = 8489.637000000001
$ws.Range("I84").Value = 8028.5713
$ws.Range("K84").Value = 80285.71299999999
$ws.Range("M84").Value = -74981.71299999999
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = ""
$ws.Range("H113").Value = 417.16666
$ws.Range("J113").Value = 493
$ws.Range("L113").Value = 1479
$ws.Range("N113").Value = -5819
$ws.Range("H124").Value = 74998.5
$ws.Range("J124").Value = 74998.5
$ws.Range("L124").Value = 74998.5
$ws.Range("N124").Value = -84818.5
$ws.Range("H126").Value = 1678.5625
$ws.Range("I126").Value = 1428.5555
$ws.Range("K126").Value = 4285.666499999999
$ws.Range("M126").Value = -1815.666499999999
$ws.Range("H139").Value = 79999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 79999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 79999
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = -90279
